# "Generate Report for Handback"
#
# Refresh the handback timestamps for the 5395e6d2-...-c9560d5a... row
# (row 3) on both the "zh-cn" and "de-de" report sheets:
#   - Column E: "Correspond Handoff Datetime"
#   - Column H: "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-14 08:30:52"
$wsZh.Range("H3").Value = "2016-03-14 08:31:16"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-14 08:30:55"
$wsDe.Range("H3").Value = "2016-03-14 08:31:26"
